$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "71.213.08"
Set-TextValue "E2" "  +3.57%  "
# Row 3
Set-TextValue "D3" "2.632.80"
Set-TextValue "E3" "  +3.86%  "
# Row 4
Set-TextValue "E4" "  -0.04%  "
# Row 5
Set-TextValue "D5" "606.78"
Set-TextValue "E5" "  +2.09%  "
# Row 6
Set-TextValue "D6" "181.11"
Set-TextValue "E6" "  +1.79%  "
# Row 7
Set-TextValue "E7" "  -0.10%  "
# Row 8
Set-TextValue "E8" "  +1.04%  "
# Row 9
Set-TextValue "D9" "2.632.10"
Set-TextValue "E9" "  +3.82%  "
# Row 10
Set-TextValue "D10" "0.166"
Set-TextValue "E10" "  +13.77%  "
# Row 11
Set-TextValue "E11" "  +0.24%  "
# Row 12
Set-TextValue "D12" "0.347"
Set-TextValue "E12" "  +2.16%  "
# Row 13
Set-TextValue "D13" "5.02"
Set-TextValue "E13" "  +0.26%  "
# Row 14
Set-TextValue "E14" "  +9.57%  "
# Row 15
Set-TextValue "D15" "3.074.91"
# Row 16
Set-TextValue "D16" "26.69"
Set-TextValue "E16" "  +2.06%  "
# Row 17
Set-TextValue "D17" "71.119.54"
Set-TextValue "E17" "  +3.54%  "
# Row 18
Set-TextValue "D18" "2.624.26"
Set-TextValue "E18" "  +4.38%  "
# Row 19
Set-TextValue "D19" "383.91"
# Row 20
Set-TextValue "D20" "7.94"
Set-TextValue "E20" "  +5.53%  "
# Row 21
Set-TextValue "D21" "11.52"
Set-TextValue "E21" "  +3.52%  "
# Row 22
Set-TextValue "E22" "  -1.88%  "
# Row 23
Set-TextValue "B23" "NEARProtocol"
Set-TextValue "C23" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D23" "4.47"
Set-TextValue "E23" "  +5.46%  "
# Row 24
Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "72.16"
Set-TextValue "E24" "  +1.28%  "
# Row 25
Set-TextValue "E25" "  +0.03%  "
# Row 26
Set-TextValue "E26" "  +10.69%  "
# Row 27
Set-TextValue "D27" "9.65"
Set-TextValue "E27" "  +6.51%  "
# Row 28
Set-TextValue "D28" "2.762.25"
Set-TextValue "E28" "  +4.67%  "
# Row 29
Set-TextValue "E29" "  -0.40%  "
# Row 30
Set-TextValue "D30" "0.0₃0965"
Set-TextValue "E30" "  +6.99%  "
# Row 31
Set-TextValue "D31" "541.86"
Set-TextValue "E31" "  +5.05%  "
# Row 32
Set-TextValue "D32" "8.06"
Set-TextValue "E32" "  +2.97%  "
# Row 33
Set-TextValue "D33" "1.33"
Set-TextValue "E33" "  +5.00%  "
# Row 34
Set-TextValue "D34" "1.84"
Set-TextValue "E34" "  +3.14%  "
# Row 35
Set-TextValue "E35" "  -0.12%  "
# Row 36
Set-TextValue "D36" "165.84"
Set-TextValue "E36" "  +1.01%  "
# Row 37
Set-TextValue "E37" "  -1.77%  "
# Row 38
Set-TextValue "D38" "19.24"
Set-TextValue "E38" "  +4.24%  "
# Row 39
Set-TextValue "D39" "1.89"
Set-TextValue "E39" "  +7.04%  "
# Row 42
Set-TextValue "E42" "  +8.63%  "
# Row 43
Set-TextValue "E43" "  +0.05%  "
# Row 44
Set-TextValue "D44" "5.05"
Set-TextValue "E44" "  +3.86%  "
# Row 45
Set-TextValue "D45" "0.331"
Set-TextValue "E45" "  +1.07%  "
# Row 46
Set-TextValue "D46" "39.98"
Set-TextValue "E46" "  +2.34%  "
# Row 47
Set-TextValue "D47" "154.28"
Set-TextValue "E47" "  +0.83%  "
# Row 48
Set-TextValue "D48" "3.64"
Set-TextValue "E48" "  +1.74%  "
# Row 49
Set-TextValue "E49" "  +5.07%  "
# Row 50
Set-TextValue "D50" "0.533"
Set-TextValue "E50" "  +2.24%  "
# Row 51
Set-TextValue "D51" "0.0₆0264"
Set-TextValue "E51" "  +0.91%  "
